$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 12.23590000000001
$ws.Range("E7").Value = 11.8924
$ws.Range("D8").Value = -8.733599999999992
$ws.Range("C12").Value = -14.58000000000002
$ws.Range("D12").Value = -8.149900000000004
$ws.Range("D14").Value = -8.752599999999999
$ws.Range("E19").Value = 12.8436
$ws.Range("E21").Value = 12.8105
$ws.Range("D22").Value = -8.001499999999997
$ws.Range("E24").Value = 12.82459999999999

$wb.Save()
